$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("2:41").Delete()

$ws.Range("A2").Value = "param_TimeStep_starting_index"
$ws.Range("B2").Value = 10
$ws.Range("A3").Value = "param_demand1_op_cost_starting_index"
$ws.Range("B3").Value = 0
$ws.Range("A4").Value = "param_demand1_inv_cost_starting_index"
$ws.Range("B4").Value = 0
$ws.Range("A5").Value = "param_net1_buy_electric_starting_index"
$ws.Range("B5").Value = 128.6777790917907
$ws.Range("A6").Value = "param_P_from_net1_starting_index"
$ws.Range("B6").Value = 273.7825087059377
$ws.Range("A7").Value = "param_net1_sell_electric_starting_index"
$ws.Range("B7").Value = 0
$ws.Range("A8").Value = "param_Q_from_net1_starting_index"
$ws.Range("B8").Value = 510.9182330044023
$ws.Range("A9").Value = "param_net1_buy_thermal_starting_index"
$ws.Range("B9").Value = 199.2581108717169
$ws.Range("A10").Value = "param_P_net1_demand1_starting_index"
$ws.Range("B10").Value = 273.7825087059377
$ws.Range("A11").Value = "param_net1_sell_thermal_starting_index"
$ws.Range("B11").Value = 0
$ws.Range("A12").Value = "param_Q_to_net1_starting_index"
$ws.Range("B12").Value = 0
$ws.Range("A13").Value = "param_P_to_net1_starting_index"
$ws.Range("B13").Value = 0
$ws.Range("A14").Value = "param_net1_emissions_starting_index"
$ws.Range("B14").Value = 275.9385807963816
$ws.Range("A15").Value = "param_net1_inv_cost_starting_index"
$ws.Range("B15").Value = 0
$ws.Range("A16").Value = "param_Q_net1_demand1_starting_index"
$ws.Range("B16").Value = 510.9182330044023
$ws.Range("A17").Value = "param_P_net1_bat1_starting_index"
$ws.Range("B17").Value = 0
$ws.Range("A18").Value = "param_pv1_op_cost_starting_index"
$ws.Range("B18").Value = 1
$ws.Range("A19").Value = "param_P_from_pv1_starting_index"
$ws.Range("B19").Value = 0.3333336
$ws.Range("A20").Value = "param_P_pv1_net1_starting_index"
$ws.Range("B20").Value = 0
$ws.Range("A21").Value = "param_pv1_inv_cost_starting_index"
$ws.Range("B21").Value = 0
$ws.Range("A22").Value = "param_pv1_emissions_starting_index"
$ws.Range("B22").Value = 0.1666668
$ws.Range("A23").Value = "param_P_pv1_demand1_starting_index"
$ws.Range("B23").Value = 0.3333336
$ws.Range("A24").Value = "param_P_pv1_bat1_starting_index"
$ws.Range("B24").Value = 0
$ws.Range("A25").Value = "param_P_from_bat1_starting_index"
$ws.Range("B25").Value = 66.49631303033051
$ws.Range("A26").Value = "param_bat1_emissions_starting_index"
$ws.Range("B26").Value = 0.184711980639807
$ws.Range("A27").Value = "param_bat1_cumulated_aging_starting_index"
$ws.Range("B27").Value = [double]"4.989292588558339E-05"
$ws.Range("A28").Value = "param_bat1_K_dis_starting_index"
$ws.Range("B28").Value = 1
$ws.Range("A29").Value = "param_P_bat1_net1_starting_index"
$ws.Range("B29").Value = 0
$ws.Range("A30").Value = "param_bat1_K_ch_starting_index"
$ws.Range("B30").Value = 0
$ws.Range("A31").Value = "param_bat1_SOC_max_starting_index"
$ws.Range("B31").Value = 0.9999501070741145
$ws.Range("A32").Value = "param_bat1_inv_cost_starting_index"
$ws.Range("B32").Value = 0
$ws.Range("A33").Value = "param_P_bat1_demand1_starting_index"
$ws.Range("B33").Value = 66.49631303033051
$ws.Range("A34").Value = "param_bat1_SOC_starting_index"
$ws.Range("B34").Value = 0.3
$ws.Range("A35").Value = "param_bat1_op_cost_starting_index"
$ws.Range("B35").Value = 1
$ws.Range("A36").Value = "param_bat1_integer_starting_index"
$ws.Range("B36").Value = -0
$ws.Range("A37").Value = "param_P_to_bat1_starting_index"
$ws.Range("B37").Value = 0
$ws.Range("A38").Value = "param_total_operation_cost_starting_index"
$ws.Range("B38").Value = 2
$ws.Range("A39").Value = "param_total_buy_starting_index"
$ws.Range("B39").Value = 327.9358899635076
$ws.Range("A40").Value = "param_total_sell_starting_index"
$ws.Range("B40").Value = 0
$ws.Range("A41").Value = "param_total_emissions_starting_index"
$ws.Range("B41").Value = 276.2899595770214
